$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add row 5 with the new trade's data
$ws.Cells.Item(5, 1).Value = 9928.61
$ws.Cells.Item(5, 2).Value = 10044.120000000001
$ws.Cells.Item(5, 3).Value = 19.170000000000002
$ws.Cells.Item(5, 4).Value = 19.39
$ws.Cells.Item(5, 5).Value = $true
$ws.Cells.Item(5, 6).Value = 1.1499999999999999
$ws.Cells.Item(5, 7).Value = 42609.505844907406
$ws.Cells.Item(5, 8).Value = $false

# Copy the date/time number format from the cell above (G4) so G5 reuses the
# existing style instead of Excel creating a brand new numFmt entry.
$ws.Range("G4").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
